$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Add blue fill/style to the data rows (A5:J15)
$rng = $ws.Range("A5:J15")
$rng.Interior.Color = 0xCCA329
$rng.Font.Name = "Arial"
$rng.Font.Size = 11
$rng.Borders.LineStyle = 1
$rng.Borders.Color = 0x000000
$rng.HorizontalAlignment = -4108
$rng.VerticalAlignment = -4108
$rng.WrapText = $true

# 2. B19 becomes boolean FALSE instead of the " " shared string
$ws.Range("B19").Value = $false

# 3. Fix FLOOR() calls that had a stray extra argument
$ws.Range("B22").Formula = '=FLOOR(F17/8,1)&"."&FLOOR(MOD(F17,8),1)&"."&(MOD(F17,8)-FLOOR(MOD(F17,8),1))*60'
$ws.Range("B23").Formula = '=FLOOR(H19,1)&"."&(H19-FLOOR(H19,1))*8&".0"'
$ws.Range("B24").Formula = '=FLOOR(I19,1)&"."&(I19-FLOOR(I19,1))*8&".0"'
$ws.Range("B27").Formula = '=FLOOR(K27/8,1)&"."&FLOOR(MOD(K27,8),1)&"."&(MOD(K27,8)-FLOOR(MOD(K27,8),1))*60'
